$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '25.901.24' }
    @{ Cell = 'E2'; Value = '  -0.54%  ' }
    @{ Cell = 'D3'; Value = '1.641.78' }
    @{ Cell = 'E3'; Value = '  +0.10%  ' }
    @{ Cell = 'D4'; Value = '''1.004' }
    @{ Cell = 'E4'; Value = '  -0.45%  ' }
    @{ Cell = 'D5'; Value = '''215.20' }
    @{ Cell = 'E5'; Value = '  -0.01%  ' }
    @{ Cell = 'D6'; Value = '''0.5056' }
    @{ Cell = 'E6'; Value = '  +0.13%  ' }
    @{ Cell = 'D7'; Value = '''1.005' }
    @{ Cell = 'E7'; Value = '  -0.41%  ' }
    @{ Cell = 'D8'; Value = '''0.2572' }
    @{ Cell = 'E8'; Value = '  -0.12%  ' }
    @{ Cell = 'D9'; Value = '''0.06395' }
    @{ Cell = 'E9'; Value = '  -0.64%  ' }
    @{ Cell = 'D10'; Value = '''19.60' }
    @{ Cell = 'E10'; Value = '  +0.77%  ' }
    @{ Cell = 'D11'; Value = '''0.07794' }
    @{ Cell = 'E11'; Value = '  +0.78%  ' }
    @{ Cell = 'B12'; Value = 'Polkadot' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot' }
    @{ Cell = 'D12'; Value = '''4.280' }
    @{ Cell = 'E12'; Value = '  +0.84%  ' }
    @{ Cell = 'B13'; Value = 'WrappedEther' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth' }
    @{ Cell = 'D13'; Value = '1.650.65' }
    @{ Cell = 'E13'; Value = '  +0.54%  ' }
    @{ Cell = 'D14'; Value = '''0.5431' }
    @{ Cell = 'E14'; Value = '  -0.40%  ' }
    @{ Cell = 'D15'; Value = '0.0₅7865' }
    @{ Cell = 'E15'; Value = '  -0.45%  ' }
    @{ Cell = 'D16'; Value = '''64.84' }
    @{ Cell = 'E16'; Value = '  +1.94%  ' }
    @{ Cell = 'D17'; Value = '25.946.31' }
    @{ Cell = 'E17'; Value = '  -0.30%  ' }
    @{ Cell = 'D19'; Value = '''198.11' }
    @{ Cell = 'E19'; Value = '  -2.94%  ' }
    @{ Cell = 'D20'; Value = '''4.390' }
    @{ Cell = 'D21'; Value = '''9.965' }
    @{ Cell = 'E21'; Value = '  -0.28%  ' }
    @{ Cell = 'D22'; Value = '''5.989' }
    @{ Cell = 'E22'; Value = '  +0.32%  ' }
    @{ Cell = 'D23'; Value = '''1.007' }
    @{ Cell = 'E23'; Value = '  -0.21%  ' }
    @{ Cell = 'D24'; Value = '''1.866' }
    @{ Cell = 'E24'; Value = '  -3.35%  ' }
    @{ Cell = 'D25'; Value = '''140.15' }
    @{ Cell = 'E25'; Value = '  -0.85%  ' }
    @{ Cell = 'E26'; Value = '  -0.89%  ' }
    @{ Cell = 'D27'; Value = '''6.854' }
    @{ Cell = 'E27'; Value = '  +1.73%  ' }
    @{ Cell = 'D28'; Value = '''15.73' }
    @{ Cell = 'E28'; Value = '  +0.02%  ' }
    @{ Cell = 'D29'; Value = '''1.244' }
    @{ Cell = 'E29'; Value = '  +0.21%  ' }
    @{ Cell = 'D30'; Value = '''0.04938' }
    @{ Cell = 'E30'; Value = '  -2.27%  ' }
    @{ Cell = 'D31'; Value = '''3.268' }
    @{ Cell = 'E31'; Value = '  +0.49%  ' }
    @{ Cell = 'D32'; Value = '''3.198' }
    @{ Cell = 'E32'; Value = '  +0.22%  ' }
    @{ Cell = 'D33'; Value = '''1.534' }
    @{ Cell = 'E33'; Value = '  -0.46%  ' }
    @{ Cell = 'D34'; Value = '''2.372' }
    @{ Cell = 'E34'; Value = '  +1.42%  ' }
    @{ Cell = 'D35'; Value = '''0.8936' }
    @{ Cell = 'E35'; Value = '  -0.09%  ' }
    @{ Cell = 'D36'; Value = '''2.607' }
    @{ Cell = 'E36'; Value = '  -0.29%  ' }
    @{ Cell = 'D37'; Value = '1.140.27' }
    @{ Cell = 'E37'; Value = '  -0.60%  ' }
    @{ Cell = 'D38'; Value = '''0.5548' }
    @{ Cell = 'E38'; Value = '  -1.51%  ' }
    @{ Cell = 'D39'; Value = '''0.01560' }
    @{ Cell = 'E39'; Value = '  -0.69%  ' }
    @{ Cell = 'D40'; Value = '''1.007' }
    @{ Cell = 'E40'; Value = '  -0.21%  ' }
    @{ Cell = 'D41'; Value = '''5.678' }
    @{ Cell = 'E41'; Value = '  +0.10%  ' }
    @{ Cell = 'D42'; Value = '''0.8207' }
    @{ Cell = 'E42'; Value = '  +1.18%  ' }
    @{ Cell = 'D43'; Value = '''99.37' }
    @{ Cell = 'E43'; Value = '  -0.42%  ' }
    @{ Cell = 'E44'; Value = '  +6.32%  ' }
    @{ Cell = 'D45'; Value = '1.778.96' }
    @{ Cell = 'E45'; Value = '  +0.05%  ' }
    @{ Cell = 'D46'; Value = '''0.4520' }
    @{ Cell = 'E46'; Value = '  -0.22%  ' }
    @{ Cell = 'D47'; Value = '''55.32' }
    @{ Cell = 'E47'; Value = '  +0.74%  ' }
    @{ Cell = 'D48'; Value = '''1.005' }
    @{ Cell = 'E48'; Value = '  -0.57%  ' }
    @{ Cell = 'D49'; Value = '''0.05052' }
    @{ Cell = 'E49'; Value = '  +0.33%  ' }
    @{ Cell = 'E50'; Value = '  -0.08%  ' }
    @{ Cell = 'D51'; Value = '''0.09510' }
    @{ Cell = 'E51'; Value = '  +2.05%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}

$wb.Save()